$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Microsoft Layoff Stories"
$ws.Range("A6").Select()
